$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 9261025
$ws.Range("I18").Value = 13889388
$ws.Range("K18").Value = 13889388
$ws.Range("M18").Value = -13889104
$ws.Range("H76").Value = 2969.6667
$ws.Range("I76").Value = 2892.3333
$ws.Range("K76").Value = 2892.3333
$ws.Range("M76").Value = -2577.3333
$ws.Range("H79").Value = 2969.6667
$ws.Range("I79").Value = 2892.3333
$ws.Range("K79").Value = 2892.3333
$ws.Range("M79").Value = -1800.3333
$ws.Range("H112").Value = 15626939
$ws.Range("J112").Value = 22728638
$ws.Range("L112").Value = 68185914
$ws.Range("N112").Value = -68188130
$ws.Range("H121").Value = 1267.6154
$ws.Range("J121").Value = 1339.9166
$ws.Range("L121").Value = 4019.7498
$ws.Range("N121").Value = -7513.7498
$ws.Range("H137").Value = 3768.7805
$ws.Range("I137").Value = 4085.037
$ws.Range("J137").Value = 3158.8572
$ws.Range("K137").Value = 12255.111
$ws.Range("L137").Value = 9476.571599999999
$ws.Range("M137").Value = -9705.110999999999
$ws.Range("N137").Value = -14576.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6228.71
$ws.Range("I32").Value = 5714.4434
$ws.Range("K32").Value = 5714.4434
$ws.Range("M32").Value = -5427.4434
$ws.Range("H45").Value = 1620.8334
$ws.Range("I45").Value = 1256.5217
$ws.Range("J45").Value = 10000
$ws.Range("K45").Value = 1256.5217
$ws.Range("L45").Value = 10000
$ws.Range("M45").Value = -879.5217
$ws.Range("N45").Value = -10754
$ws.Range("H88").Value = 4669
$ws.Range("J88").Value = 6003.5
$ws.Range("L88").Value = 6003.5
$ws.Range("N88").Value = -6815.5
$ws.Range("H91").Value = 4669
$ws.Range("J91").Value = 6003.5
$ws.Range("L91").Value = 6003.5
$ws.Range("N91").Value = -8811.5
$ws.Range("H130").Value = 46886.332
$ws.Range("J130").Value = 46886.332
$ws.Range("L130").Value = 46886.332
$ws.Range("N130").Value = -56926.332
$ws.Range("H137").Value = 29531.666
$ws.Range("J137").Value = 29531.666
$ws.Range("L137").Value = 29531.666
$ws.Range("N137").Value = -39731.666
$ws.Range("H139").Value = 29642.857
$ws.Range("J139").Value = 29642.857
$ws.Range("L139").Value = 29642.857
$ws.Range("N139").Value = -39922.857
$ws.Range("H140").Value = 29642.857
$ws.Range("J140").Value = 29642.857
$ws.Range("L140").Value = 29642.857
$ws.Range("N140").Value = -40002.857
$ws.Range("H141").Value = 30000
$ws.Range("J141").Value = 30000
$ws.Range("L141").Value = 30000
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 100009
$ws.Range("J16").Value = 100009
$ws.Range("L16").Value = 100009
$ws.Range("N16").Value = -100349
$ws.Range("H23").Value = 67220.14
$ws.Range("J23").Value = 84008.2
$ws.Range("L23").Value = 84008.2
$ws.Range("N23").Value = -84574.2
$ws.Range("H68").Value = 20268
$ws.Range("I68").Value = 20268
$ws.Range("K68").Value = 20268
$ws.Range("M68").Value = -19457
$ws.Range("H71").Value = 20268
$ws.Range("I71").Value = 20268
$ws.Range("K71").Value = 60804
$ws.Range("M71").Value = -56748
$ws.Range("H86").Value = 255000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 255000
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = 255000
$ws.Range("N86").Value = -257246
$ws.Range("H89").Value = 255000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 255000
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = 1275000
$ws.Range("N89").Value = -1286232
$ws.Range("H99").Value = 2719.8948
$ws.Range("I99").Value = 2451.8667
$ws.Range("J99").Value = 3725
$ws.Range("K99").Value = 2451.8667
$ws.Range("L99").Value = 3725
$ws.Range("M99").Value = -953.8667
$ws.Range("N99").Value = -6721
$ws.Range("H107").Value = 1932.1578
$ws.Range("I107").Value = 1491.5
$ws.Range("K107").Value = 1491.5
$ws.Range("M107").Value = 428.5
$ws.Range("H130").Value = 27500
$ws.Range("J130").Value = 27500
$ws.Range("L130").Value = 27500
$ws.Range("N130").Value = -37540
$ws.Range("H134").Value = 6389.84
$ws.Range("I134").Value = 6606.952
$ws.Range("J134").Value = 5250
$ws.Range("K134").Value = 19820.856
$ws.Range("L134").Value = 15750
$ws.Range("M134").Value = -17285.856
$ws.Range("N134").Value = -20820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 2188.0908
$ws.Range("I138").Value = 1218.4286
$ws.Range("J138").Value = 3885
$ws.Range("K138").Value = 3655.2858
$ws.Range("L138").Value = 11655
$ws.Range("M138").Value = 1484.7142
$ws.Range("N138").Value = -21935

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 30000
$ws.Range("J130").Value = 30000
$ws.Range("L130").Value = 30000
$ws.Range("N130").Value = -40040
$ws.Range("H132").Value = 5416.037
$ws.Range("I132").Value = 7969.5835
$ws.Range("J132").Value = 3373.2
$ws.Range("K132").Value = 23908.7505
$ws.Range("L132").Value = 10119.6
$ws.Range("M132").Value = -21378.7505
$ws.Range("N132").Value = -15179.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 627.1
$ws.Range("I9").Value = 201.66667
$ws.Range("J9").Value = 1265.25
$ws.Range("K9").Value = 201.66667
$ws.Range("L9").Value = 1265.25
$ws.Range("M9").Value = 22.33332999999999
$ws.Range("N9").Value = -1713.25
$ws.Range("H94").Value = 31500
$ws.Range("J94").Value = 31500
$ws.Range("L94").Value = 31500
$ws.Range("N94").Value = -32852
$ws.Range("H138").Value = 30237.5
$ws.Range("I138").Value = 10000
$ws.Range("J138").Value = 31586.666
$ws.Range("K138").Value = 10000
$ws.Range("L138").Value = 31586.666
$ws.Range("M138").Value = -4860
$ws.Range("N138").Value = -41866.666
$ws.Range("H140").Value = 29519.637
$ws.Range("J140").Value = 29519.637
$ws.Range("L140").Value = 29519.637
$ws.Range("N140").Value = -39879.637
$ws.Range("H141").Value = 29818.182
$ws.Range("J141").Value = 29818.182
$ws.Range("L141").Value = 29818.182
$ws.Range("N141").Value = -40178.182

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 29899.75
$ws.Range("J92").Value = 29899.75
$ws.Range("L92").Value = 29899.75
$ws.Range("N92").Value = -34891.75
$ws.Range("H113").Value = 949.7692
$ws.Range("I113").Value = 175.66667
$ws.Range("K113").Value = 527.00001
$ws.Range("M113").Value = 1642.99999
$ws.Range("H133").Value = 36666.668
$ws.Range("J133").Value = 36666.668
$ws.Range("L133").Value = 36666.668
$ws.Range("N133").Value = -46786.668
$ws.Range("H135").Value = 50153.75
$ws.Range("J135").Value = 50153.75
$ws.Range("L135").Value = 50153.75
$ws.Range("N135").Value = -60293.75
$ws.Range("H137").Value = 34767.145
$ws.Range("J137").Value = 34767.145
$ws.Range("L137").Value = 34767.145
$ws.Range("N137").Value = -44967.145
$ws.Range("H140").Value = 29581.818
$ws.Range("J140").Value = 29581.818
$ws.Range("L140").Value = 29581.818
$ws.Range("N140").Value = -39941.818
$ws.Range("H141").Value = 28687.5
$ws.Range("J141").Value = 28687.5
$ws.Range("L141").Value = 28687.5
$ws.Range("N141").Value = -39047.5
